# "Generate Report for Archive"
#
# The localization status report moved from "Ready for handoff" to
# "In Translation" for the two source files tracked in this workbook.
# Update the status text everywhere it appears (the Overview rollup
# sheet plus each per-locale sheet), then shrink the now-narrower
# status columns to fit the new (shorter) text, mirroring Excel's
# column-width auto-adjustment that accompanied the original edit.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# --- Overview sheet: status is duplicated per-locale in columns E (zh-cn)
# and F (de-de) for both tracked files (rows 2-3) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# --- Per-locale sheets: status lives in column C ("Status") ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# --- Shrink the status columns to fit the shorter text ---
# (New text is shorter than "Ready for handoff", so the column narrows.)
# Use numeric column indices: Overview E=5, F=6; per-locale sheets C=3.
$newColumnWidth = 12.5

$overview.Columns.Item(5).ColumnWidth = $newColumnWidth
$overview.Columns.Item(6).ColumnWidth = $newColumnWidth
$zhcn.Columns.Item(3).ColumnWidth = $newColumnWidth
$dede.Columns.Item(3).ColumnWidth = $newColumnWidth

Write-Output "Status updated to '$newStatus' and status columns resized."
